{"js": "// Locate the paragraph \"6. Refatorar o c\u00f3digo.\" and:\n//  1) insert a new paragraph right after it with text \"8. Refatorar o c\u00f3digo.\"\n//     (inherits the same run formatting: Helvetica, sz 58, from the source\n//     paragraph mark)\n//  2) change the original paragraph's own text to\n//     \"7. Mostrar o aviso de ganhou ou perdeu.\" (keeps its formatting)\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst OLD_TEXT = \"6. Refatorar o c\u00f3digo.\";\nconst NEW_TEXT_7 = \"7. Mostrar o aviso de ganhou ou perdeu.\";\nconst NEW_TEXT_8 = \"8. Refatorar o c\u00f3digo.\";\n\nlet target = null;\nfor (const p of paragraphs.items) {\n  if (p.text === OLD_TEXT) {\n    target = p;\n    break;\n  }\n}\n\nif (!target) {\n  throw new Error(`Could not find paragraph '${OLD_TEXT}'`);\n}\n\n// Insert the new \"8. ...\" paragraph directly after the target while its\n// text/formatting are still intact, so the new paragraph's run inherits the\n// same rPr (Helvetica, sz 58) as the source paragraph mark.\ntarget.insertParagraph(NEW_TEXT_8, \"After\");\n\n// Now replace the original paragraph's text in place, preserving formatting.\ntarget.getRange(\"Whole\").insertText(NEW_TEXT_7, \"Replace\");\n\nawait context.sync();\n", "ps1": "# Locate the paragraph \"6. Refatorar o c\u00f3digo.\" and:\n#  1) insert a new paragraph right after it with text \"8. Refatorar o c\u00f3digo.\"\n#     (inherits the same run formatting: Helvetica, sz 58)\n#  2) change the original paragraph's own text to\n#     \"7. Mostrar o aviso de ganhou ou perdeu.\" (keeps its formatting)\n\n$d = $word.ActiveDocument\n\n$OLD_TEXT = \"6. Refatorar o c\u00f3digo.\"\n$NEW_TEXT_7 = \"7. Mostrar o aviso de ganhou ou perdeu.\"\n$NEW_TEXT_8 = \"8. Refatorar o c\u00f3digo.\"\n\n# Find the 1-based paragraph index of the target paragraph.\n$targetIndex = -1\n$i = 0\nforeach ($p in $d.Paragraphs) {\n    $i = $i + 1\n    $txt = $p.Range.Text -replace \"[\\r\\a\\v]+$\", \"\"\n    if ($txt -eq $OLD_TEXT) {\n        $targetIndex = $i\n        break\n    }\n}\n\nif ($targetIndex -eq -1) {\n    throw \"Could not find paragraph '$OLD_TEXT'\"\n}\n\n# Insert a new paragraph mark right after the target paragraph while its\n# text/formatting are still intact, so the new paragraph's run inherits the\n# same rPr (Helvetica, sz 58) as the source paragraph mark.\n$targetRange = $d.Paragraphs.Item($targetIndex).Range\n$targetRange.InsertParagraphAfter()\n\n# The freshly inserted paragraph is now immediately after the target.\n$newParaRange = $d.Paragraphs.Item($targetIndex + 1).Range\n$newParaRange.Text = $NEW_TEXT_8\n\n# Update the original paragraph's text in place, preserving its formatting.\n$origRange = $d.Paragraphs.Item($targetIndex).Range\n$origRange.Text = $NEW_TEXT_7\n"}
